$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 20.97917533333333
$ws.Range("H2").Value = 62.93752600000001
$ws.Range("I2").Value = 0.2451892257562263
$ws.Range("J2").Value = 0.2451892257562263
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 9.266076
$ws.Range("N2").Value = 27.798228
$ws.Range("O2").Value = 0.5506394579555814
$ws.Range("P2").Value = 0.5506394579555816
$ws.Range("Q2").Value = 194.394633055992
$ws.Range("R2").Value = 1749.551697503928
$ws.Range("S2").Value = 0.1350108623669571
$ws.Range("T2").Value = 0.1350108623669571

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 20.97917533333333
$ws.Range("H3").Value = 62.93752600000001
$ws.Range("I3").Value = 0.2451892257562263
$ws.Range("J3").Value = 0.2451892257562263
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 2.181666333333333
$ws.Range("N3").Value = 6.544999
$ws.Range("O3").Value = 0.1296462026888844
$ws.Range("P3").Value = 0.1296462026888844
$ws.Range("Q3").Value = 45.76956052583044
$ws.Range("R3").Value = 411.926044732474
$ws.Range("S3").Value = 0.03178785205952236
$ws.Range("T3").Value = 0.03178785205952236

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 20.97917533333333
$ws.Range("H4").Value = 62.93752600000001
$ws.Range("I4").Value = 0.2451892257562263
$ws.Range("J4").Value = 0.2451892257562263
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 5.380103666666667
$ws.Range("N4").Value = 16.140311
$ws.Range("O4").Value = 0.3197143393555341
$ws.Range("P4").Value = 0.3197143393555341
$ws.Range("Q4").Value = 112.8701381345096
$ws.Range("R4").Value = 1015.831243210586
$ws.Range("S4").Value = 0.07839051132974678
$ws.Range("T4").Value = 0.07839051132974678

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 59.05285266666667
$ws.Range("H5").Value = 177.158558
$ws.Range("I5").Value = 0.6901664624076501
$ws.Range("J5").Value = 0.6901664624076501
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 9.266076
$ws.Range("N5").Value = 27.798228
$ws.Range("O5").Value = 0.5506394579555814
$ws.Range("P5").Value = 0.5506394579555816
$ws.Range("Q5").Value = 547.188220826136
$ws.Range("R5").Value = 4924.693987435225
$ws.Range("S5").Value = 0.3800328867592697
$ws.Range("T5").Value = 0.3800328867592697

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 59.05285266666667
$ws.Range("H6").Value = 177.158558
$ws.Range("I6").Value = 0.6901664624076501
$ws.Range("J6").Value = 0.6901664624076501
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 2.181666333333333
$ws.Range("N6").Value = 6.544999
$ws.Range("O6").Value = 0.1296462026888844
$ws.Range("P6").Value = 0.1296462026888844
$ws.Range("Q6").Value = 128.8336205501602
$ws.Range("R6").Value = 1159.502584951442
$ws.Range("S6").Value = 0.08947746107437254
$ws.Range("T6").Value = 0.08947746107437254

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 59.05285266666667
$ws.Range("H7").Value = 177.158558
$ws.Range("I7").Value = 0.6901664624076501
$ws.Range("J7").Value = 0.6901664624076501
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 5.380103666666667
$ws.Range("N7").Value = 16.140311
$ws.Range("O7").Value = 0.3197143393555341
$ws.Range("P7").Value = 0.3197143393555341
$ws.Range("Q7").Value = 317.7104691590598
$ws.Range("R7").Value = 2859.394222431538
$ws.Range("S7").Value = 0.2206561145740079
$ws.Range("T7").Value = 0.2206561145740079

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 5.531174333333333
$ws.Range("H8").Value = 16.593523
$ws.Range("I8").Value = 0.06464431183612354
$ws.Range("J8").Value = 0.06464431183612354
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 9.266076
$ws.Range("N8").Value = 27.798228
$ws.Range("O8").Value = 0.5506394579555814
$ws.Range("P8").Value = 0.5506394579555816
$ws.Range("Q8").Value = 51.252281741916
$ws.Range("R8").Value = 461.2705356772441
$ws.Range("S8").Value = 0.03559570882935464
$ws.Range("T8").Value = 0.03559570882935465

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 5.531174333333333
$ws.Range("H9").Value = 16.593523
$ws.Range("I9").Value = 0.06464431183612354
$ws.Range("J9").Value = 0.06464431183612354
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 2.181666333333333
$ws.Range("N9").Value = 6.544999
$ws.Range("O9").Value = 0.1296462026888844
$ws.Range("P9").Value = 0.1296462026888844
$ws.Range("Q9").Value = 12.06717682683078
$ws.Range("R9").Value = 108.604591441477
$ws.Range("S9").Value = 0.008380889554989524
$ws.Range("T9").Value = 0.008380889554989524

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 5.531174333333333
$ws.Range("H10").Value = 16.593523
$ws.Range("I10").Value = 0.06464431183612354
$ws.Range("J10").Value = 0.06464431183612354
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 5.380103666666667
$ws.Range("N10").Value = 16.140311
$ws.Range("O10").Value = 0.3197143393555341
$ws.Range("P10").Value = 0.3197143393555341
$ws.Range("Q10").Value = 29.75829131173922
$ws.Range("R10").Value = 267.824621805653
$ws.Range("S10").Value = 0.02066771345177937
$ws.Range("T10").Value = 0.02066771345177937

